# Add two new payslip history rows (employee #4 and #5) to the employee_data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: employee_id is numeric
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "RyanKho"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5

# Row 6: employee_id stored as text "5" (matches source data quirk)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "Thim"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 6
